{"js": "const body = context.document.body;\nconst paras = body.paragraphs;\nparas.load(\"items/text\");\nawait context.sync();\n\n// Find the bullet-list paragraph whose text is exactly \"[Antorus]\"\nlet target = null;\nfor (const p of paras.items) {\n  if (p.text === \"[Antorus]\") {\n    target = p;\n    break;\n  }\n}\nif (!target) {\n  throw new Error(\"Could not find paragraph with text '[Antorus]'\");\n}\n\n// Insert the three new bullet-list items right after it, preserving order.\nlet anchor = target;\nconst newItems = [\"[Megas]\", \"[Archaon]\", \"[Ho]\"];\nfor (const txt of newItems) {\n  anchor = anchor.insertParagraph(txt, \"After\");\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the bullet-list item \"[Antorus]\" \u2014 the new items are inserted\n# directly after it, keeping their place in the same bulleted list.\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$found = $rng.Find.Execute(\"[Antorus]\")\nif (-not $found) {\n  throw \"Could not find paragraph with text '[Antorus]'\"\n}\n\n$anchorIndex = $rng.Paragraphs.Item(1).Index\n\n$newItems = @(\"[Megas]\", \"[Archaon]\", \"[Ho]\")\nforeach ($txt in $newItems) {\n  $p = $d.Paragraphs.Item($anchorIndex)\n  $p.Range.InsertParagraphAfter()\n  $anchorIndex = $anchorIndex + 1\n  $newPara = $d.Paragraphs.Item($anchorIndex)\n  $newPara.Range.Text = $txt\n}\n"}
